# Add a new "2022-Q3" sheet (before "2022-Q2"), populated with the
# 2022-Q3 fund holdings, and record the new quarter in the "总计" summary
# sheet (inserting a row at the top of the data and renumbering the
# existing rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q3" worksheet by duplicating "2022-Q2" so it
#    inherits the same header row / column-A styling, then rename it
#    and drop it into place right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The duplicated sheet only has 5 data rows (rows 2-6); 2022-Q3 needs
# 10, so extend the row block (copying just the used columns, never a
# whole row, to avoid touching all 16384 columns).
$q3.Range("A6:H6").Copy()
for ($i = 7; $i -le 11; $i++) {
    $q3.Range("A" + $i + ":H" + $i).PasteSpecial(-4122)
}

# Helper: force a value to be stored as text (keeps leading zeros /
# trailing zeros intact) without leaving the quote-prefix style behind.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.ClearFormats()
}

$q3rows = @(
    @{A=0; B="000968"; C="广发中证养老产业指数A";         D="10.84"; E="93.98"; F="1.81"; G="0.1962"; H=1},
    @{A=1; B="001243"; C="博时中证淘金大数据100指数I";     D="2.25";  E="93.23"; F="0.91"; G="0.0205"; H=8},
    @{A=2; B="002982"; C="广发中证养老产业指数C";         D="0.92";  E="93.98"; F="1.81"; G="0.0167"; H=1},
    @{A=3; B="001242"; C="博时中证淘金大数据100指数A";     D="1.59";  E="93.23"; F="0.91"; G="0.0145"; H=8},
    @{A=4; B="013878"; C="圆信永丰中证500指数增强A";       D="0.95";  E="92.59"; F="1.48"; G="0.0141"; H=7},
    @{A=5; B="516560"; C="华宝养老ETF";                   D="0.74";  E="98.01"; F="1.88"; G="0.0139"; H=1},
    @{A=6; B="008124"; C="中邮中证500指数增强C";           D="0.70";  E="93.79"; F="1.54"; G="0.0108"; H=8},
    @{A=7; B="590007"; C="中邮中证500指数增强A";           D="0.27";  E="93.79"; F="1.54"; G="0.0042"; H=8},
    @{A=8; B="013879"; C="圆信永丰中证500指数增强C";       D="0.09";  E="92.59"; F="1.48"; G="0.0013"; H=7},
    @{A=9; B="015245"; C="南华丰汇混合";                   D="0.09";  E="86.53"; F="1.06"; G="0.0010"; H=6}
)

for ($i = 0; $i -lt $q3rows.Count; $i++) {
    $row = 2 + $i
    $data = $q3rows[$i]

    $q3.Range("A" + $row).Value = $data.A
    Set-TextValue $q3.Range("B" + $row) $data.B
    Set-TextValue $q3.Range("C" + $row) $data.C
    Set-TextValue $q3.Range("D" + $row) $data.D
    Set-TextValue $q3.Range("E" + $row) $data.E
    Set-TextValue $q3.Range("F" + $row) $data.F
    Set-TextValue $q3.Range("G" + $row) $data.G
    $q3.Range("H" + $row).Value = $data.H
}

# ---------------------------------------------------------------------
# 2. Update the "总计" (overview) sheet: insert a new row for 2022-Q3
#    above the existing data and renumber the index column.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Re-apply the style used by the other index cells (column A) by
# copying the format from the row just below.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 10
$total.Range("D2").Value = 0.29

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
